$d = $word.ActiveDocument

# --- Step 1: split paragraph 2's "hi" run into two runs: "H" and "i" ---
$p2 = $d.Paragraphs.Item(2)
$textRange = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$splitXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>H</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>i</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $textRange.InsertXML($splitXml)

# --- Step 2: add a new paragraph after it with the merge message ---
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Text = "Merged with branch2 with master"
